# 23.11.2020 MC Sales Details
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Push the existing T130 / T140 / T180 rows (A39:C41) down by one row, into
# A40:C42, then put the new D82 row into A39:C39. Only columns A:C move;
# columns F:I are untouched. Work from the bottom up so values aren't
# clobbered before they are copied.
$a41 = $ws.Range("A41").Value()
$b41 = $ws.Range("B41").Value()
$c41 = $ws.Range("C41").Value()

$a40 = $ws.Range("A40").Value()
$b40 = $ws.Range("B40").Value()
$c40 = $ws.Range("C40").Value()

$a39 = $ws.Range("A39").Value()
$b39 = $ws.Range("B39").Value()
$c39 = $ws.Range("C39").Value()

$ws.Range("A42").Value = $a41
$ws.Range("B42").Value = $b41
$ws.Range("C42").Value = $c41

$ws.Range("A41").Value = $a40
$ws.Range("B41").Value = $b40
$ws.Range("C41").Value = $c40

$ws.Range("A40").Value = $a39
$ws.Range("B40").Value = $b39
$ws.Range("C40").Value = $c39

$ws.Range("A39").Value = "D82"
$ws.Range("B39").Value = 1170
$ws.Range("C39").Value = 1250

# Update the "Last Update" label in the merged header cell F4:I4
$ws.Range("F4").Value = "Last Update(23-11-2020)"

# Update the view: clear the frozen/scrolled top-left cell and move the
# active selection to N13
$ws.Range("N13").Select() | Out-Null
